# Add 7 new data rows (23-29) to Sheet1, matching the new day (2025-05-13 / serial 45790)
# that was appended to the landscaping data log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Row, PlantType(B), PlantSize(C), Low(D), High(E), HumidityMode(H) ["v"|"f", value],
#           Pruned(I), Quadrant(J), Shade(K)
$rows = @(
    @{ R=23; B="Flowering";    C="Large";  D=61; E=66; HMode="v"; HVal=2.5;  I="No"; J=2; K="Bright" },
    @{ R=24; B="Nonflowering"; C="Medium"; D=61; E=66; HMode="v"; HVal=0.5;  I="No"; J=3; K="Bright" },
    @{ R=25; B="Nonflowering"; C="Small";  D=61; E=66; HMode="v"; HVal=0.5;  I="No"; J=3; K="Bright" },
    @{ R=26; B="Nonflowering"; C="Medium"; D=61; E=66; HMode="f"; HVal="2/6"; I="No"; J=3; K="Neutral" },
    @{ R=27; B="Nonflowering"; C="Medium"; D=61; E=66; HMode="v"; HVal=0.25; I="No"; J=3; K="Neutral" },
    @{ R=28; B="Nonflowering"; C="Large";  D=61; E=66; HMode="f"; HVal="2/3"; I="No"; J=4; K="Dark" },
    @{ R=29; B="Tree";         C="Medium"; D=61; E=66; HMode="f"; HVal="7/3"; I="No"; J=1; K="Dark" }
)

foreach ($row in $rows) {
    $r = $row.R

    # Column A: date serial, copy the formatting (numFmt) from the last existing date row (A22)
    $ws.Range("A22").Copy($ws.Range("A$r"))
    $ws.Range("A$r").Value = 45790

    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E

    # Column F: Temp_Diff = E - D
    $ws.Range("F$r").Formula = "=E$r-D$r"

    $ws.Range("G$r").Value = 2.37

    if ($row.HMode -eq "f") {
        $ws.Range("H$r").Formula = "=" + $row.HVal
    } else {
        $ws.Range("H$r").Value = $row.HVal
    }

    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K

    $ws.Range("L$r").Value = 5
    $ws.Range("M$r").Value = 0.97
    $ws.Range("N$r").Value = 64
    $ws.Range("O$r").Value = 29.87
    $ws.Range("P$r").Value = 19
    $ws.Range("Q$r").Value = 0.94
    $ws.Range("R$r").Value = 9.9
    $ws.Range("S$r").Value = 25
    $ws.Range("T$r").Value = 50
}

# Update the view: scroll so row 5 is at top, and select A30 (the empty cell below the new data)
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("A30").Select() | Out-Null
